$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column E (boolean) values to TRUE for rows 11-14, matching the diff
$ws.Range("E11").Value = $true
$ws.Range("E12").Value = $true
$ws.Range("E13").Value = $true
$ws.Range("E14").Value = $true
